# "re-ran template generation script (#576)"
#
# The cases_disposed_by_gender sheet is regenerated as cases_disposed_by_sex:
#   - gender (Male/Female/Other/Non-Binary/Unknown) is replaced by
#     biological_sex (Male Biological Sex/Female Biological Sex/Unknown Biological Sex)
#   - the data grid is rebuilt: 2 years x 12 months x 3 biological_sex values = 72 rows
#     (was 2 years x 12 months x 5 genders = 120 rows)
#   - column C is widened to fit the longer labels

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cases_disposed_by_gender")

# 1. Rename the sheet itself.
$ws.Name = "cases_disposed_by_sex"

# 2. Update the header row.
$ws.Cells.Item(1, 3).Value = "biological_sex"

# 3. Rebuild the data grid.
$years = @("2021", "2022")
$months = @("1", "2", "3", "4", "5", "6", "7", "8", "9", "10", "11", "12")
$sexes = @("Male Biological Sex", "Female Biological Sex", "Unknown Biological Sex")

$totalRows = $years.Length * $months.Length * $sexes.Length   # 72

$grid = New-Object 'object[,]' $totalRows, 3
$r = 0
foreach ($y in $years) {
    foreach ($m in $months) {
        foreach ($s in $sexes) {
            $grid[$r, 0] = $y
            $grid[$r, 1] = $m
            $grid[$r, 2] = $s
            $r = $r + 1
        }
    }
}

$firstDataRow = 2
$lastDataRow = $firstDataRow + $totalRows - 1   # 73

$dataRange = $ws.Range("A" + $firstDataRow + ":C" + $lastDataRow)

# Columns A (year) and B (month) hold values that look numeric ("2021", "1", ...)
# but must stay text, matching the workbook's existing convention (they are
# shared strings throughout, never numbers). Force text storage for the
# year/month columns, write the values, then drop back to the default
# "Normal" style so no cell keeps a lingering custom number format.
$yearMonthRange = $ws.Range("A" + $firstDataRow + ":B" + $lastDataRow)
$yearMonthRange.NumberFormat = "@"
$dataRange.Value = $grid
$yearMonthRange.Style = "Normal"

# 4. Drop the now-unused trailing rows left over from the old 120-row grid
#    (121 rows of data + header vs. the new 73).
$oldLastRow = 121
if ($oldLastRow -gt $lastDataRow) {
    $ws.Range("A" + ($lastDataRow + 1) + ":D" + $oldLastRow).EntireRow.Delete()
}

# 5. Widen column C to fit the longer "... Biological Sex" labels.
$ws.Columns.Item(3).ColumnWidth = 21.8
